$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.874076724052429
$ws.Range("B1").Value = 5.519608974456787
$ws.Range("C1").Value = 2.584070205688477
$ws.Range("D1").Value = 1.638672471046448
$ws.Range("E1").Value = 1.322680592536926
